# Trade #112 closed at 2026-02-17 16:02:51 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---- Summary sheet: refresh aggregate stats after the new trade ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1198.87
$wsSummary.Range("B4").Value = -1.14
$wsSummary.Range("B6").Value = 112
$wsSummary.Range("B8").Value = 56
$wsSummary.Range("B9").Value = 35.71

# ---- Strategy Status sheet: refresh the MarketMaking strategy row ----
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 98.87
$wsStatus.Range("D4").Value = 112
$wsStatus.Range("E4").Value = -1.14
$wsStatus.Range("F4").Value = -1.13
$wsStatus.Range("G4").Value = 35.71

# ---- Append new trade #112 to both the "All Trades" log and the
#      per-strategy "MarketMaking" log (kept in sync, identical rows) ----
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $r = 113

    $ws.Cells.Item($r, 1).Value = 112

    # Date column looks like a date string ("2026-02-17") - format the
    # cell as Text first so Excel stores it verbatim instead of coercing
    # it into a date serial number, then drop back to the Normal style.
    $dateCell = $ws.Cells.Item($r, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.Style = "Normal"

    $ws.Cells.Item($r, 3).Value = "16:02:44"
    $ws.Cells.Item($r, 4).Value = "MarketMaking"
    $ws.Cells.Item($r, 5).Value = "DOWN"
    $ws.Cells.Item($r, 6).Value = 0.21
    $ws.Cells.Item($r, 7).Value = 0.16
    $ws.Cells.Item($r, 8).Value = "CLOSED"
    $ws.Cells.Item($r, 9).Value = -23.8095
    $ws.Cells.Item($r, 10).Value = -0.05
    $ws.Cells.Item($r, 11).Value = 98.87
    $ws.Cells.Item($r, 12).Value = 0
    $ws.Cells.Item($r, 13).Value = 0
    $ws.Cells.Item($r, 14).Value = 0.6
    $ws.Cells.Item($r, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($r, 16).Value = "early_exit"
    $ws.Cells.Item($r, 17).Value = 0.14
}
